$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the date text in A6 (shared string "Ngày: 09/04/2023" -> "Ngày: 16/04/2023")
$ws.Range("A6").Value = "Ngày: 16/04/2023"

# 2) Update quantity/amount for the first item row (row 13)
$ws.Range("C13").Value = 1808
$ws.Range("F13").Value = 3254400

# 3) Update quantity/amount for the second item row (row 14)
$ws.Range("C14").Value = 2
$ws.Range("F14").Value = 3600

# 4) Remove item rows 3, 4 and 5 (old rows 15-17), shifting everything below up by 3 rows.
#    This turns the old total row (18) into row 15, etc.
$ws.Rows("15:17").Delete()

# 5) Update the total ("Cộng tiền") amount, now located at F15
$ws.Range("F15").Value = 3258000
